$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 30000
$ws.Range("J21").Value = 30000
$ws.Range("L21").Value = 30000
$ws.Range("N21").Value = -30936

$ws.Range("H23").Value = 30000
$ws.Range("J23").Value = 30000
$ws.Range("L23").Value = 30000
$ws.Range("N23").Value = -30468

$ws.Range("H43").Value = 4277.381
$ws.Range("J43").Value = 4350.4
$ws.Range("L43").Value = 4350.4
$ws.Range("N43").Value = -4488.4

$ws.Range("H64").Value = 9300
$ws.Range("I64").Value = 9900
$ws.Range("K64").Value = 9900
$ws.Range("M64").Value = -9652

$ws.Range("H67").Value = 9300
$ws.Range("I67").Value = 9900
$ws.Range("K67").Value = 9900
$ws.Range("M67").Value = -9042

$ws.Range("H74").Value = 10445.728
$ws.Range("I74").Value = 8301
$ws.Range("K74").Value = 8301
$ws.Range("M74").Value = -7365

$ws.Range("H77").Value = 10445.728
$ws.Range("I77").Value = 8301
$ws.Range("K77").Value = 41505
$ws.Range("M77").Value = -36825

$ws.Range("H138").Value = 6373.8696
$ws.Range("I138").Value = 3387.3333
$ws.Range("J138").Value = 8293.786
$ws.Range("K138").Value = 10161.9999
$ws.Range("L138").Value = 24881.358
$ws.Range("M138").Value = -5021.999899999999
$ws.Range("N138").Value = -35161.358

$ws.Range("H141").Value = 4848.6665
$ws.Range("I141").Value = 4841.1
$ws.Range("K141").Value = 14523.3
$ws.Range("M141").Value = -9343.300000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 901
$ws.Range("I4").Value = 901
$ws.Range("K4").Value = 901
$ws.Range("M4").Value = -785

$ws.Range("H45").Value = 3455.5
$ws.Range("I45").Value = 1851.8889
$ws.Range("J45").Value = 8266.333000000001
$ws.Range("K45").Value = 1851.8889
$ws.Range("L45").Value = 8266.333000000001
$ws.Range("M45").Value = -1474.8889
$ws.Range("N45").Value = -9020.333000000001

$ws.Range("H61").Value = 4921.231
$ws.Range("I61").Value = 3003.375
$ws.Range("J61").Value = 7989.8
$ws.Range("K61").Value = 3003.375
$ws.Range("L61").Value = 7989.8
$ws.Range("M61").Value = -2791.375
$ws.Range("N61").Value = -8413.799999999999

$ws.Range("H132").Value = 4309.253
$ws.Range("I132").Value = 3625.4
$ws.Range("K132").Value = 10876.2
$ws.Range("M132").Value = -8346.200000000001

$ws.Range("H136").Value = 4921.231
$ws.Range("I136").Value = 3003.375
$ws.Range("J136").Value = 7989.8
$ws.Range("K136").Value = 9010.125
$ws.Range("L136").Value = 23969.4
$ws.Range("M136").Value = -6460.125
$ws.Range("N136").Value = -29069.4

$ws.Range("H141").Value = 90000
$ws.Range("J141").Value = 90000
$ws.Range("L141").Value = 90000
$ws.Range("N141").Value = -100360

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 47192.816
$ws.Range("I20").Value = 1667
$ws.Range("K20").Value = 1667
$ws.Range("M20").Value = -1420

$ws.Range("H34").Value = 5000
$ws.Range("J34").Value = 5000
$ws.Range("L34").Value = 5000
$ws.Range("N34").Value = -5228

$ws.Range("H81").Value = 33612.25
$ws.Range("J81").Value = 33612.25
$ws.Range("L81").Value = 33612.25
$ws.Range("N81").Value = -35734.25

$ws.Range("H84").Value = 33612.25
$ws.Range("J84").Value = 33612.25
$ws.Range("L84").Value = 100836.75
$ws.Range("N84").Value = -111444.75

$ws.Range("H86").Value = 5763.222
$ws.Range("I86").Value = 5839
$ws.Range("J86").Value = 5668.5
$ws.Range("K86").Value = 5839
$ws.Range("L86").Value = 5668.5
$ws.Range("M86").Value = -4716
$ws.Range("N86").Value = -7914.5

$ws.Range("H89").Value = 5763.222
$ws.Range("I89").Value = 5839
$ws.Range("J89").Value = 5668.5
$ws.Range("K89").Value = 29195
$ws.Range("L89").Value = 28342.5
$ws.Range("M89").Value = -23579
$ws.Range("N89").Value = -39574.5

$ws.Range("H105").Value = 3835.1428
$ws.Range("I105").Value = 2914.1428
$ws.Range("J105").Value = 4756.143
$ws.Range("K105").Value = 2914.1428
$ws.Range("L105").Value = 4756.143
$ws.Range("M105").Value = -1167.1428
$ws.Range("N105").Value = -8250.143

$ws.Range("H134").Value = 4358.7144
$ws.Range("I134").Value = 3502.28
$ws.Range("J134").Value = 6499.8
$ws.Range("K134").Value = 10506.84
$ws.Range("L134").Value = 19499.4
$ws.Range("M134").Value = -7971.84
$ws.Range("N134").Value = -24569.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3750.75
$ws.Range("I58").Value = 2251.125
$ws.Range("K58").Value = 2251.125
$ws.Range("M58").Value = -2048.125

$ws.Range("H62").Value = 9658.166999999999
$ws.Range("I62").Value = 9689.799999999999
$ws.Range("K62").Value = 9689.799999999999
$ws.Range("M62").Value = -9065.799999999999

$ws.Range("H65").Value = 9658.166999999999
$ws.Range("I65").Value = 9689.799999999999
$ws.Range("K65").Value = 48449
$ws.Range("M65").Value = -45329

$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()

$ws.Range("H134").Value = 5932.7
$ws.Range("I134").Value = 4133.846
$ws.Range("K134").Value = 12401.538
$ws.Range("M134").Value = -9866.537999999999

$ws.Range("H136").Value = 3750.75
$ws.Range("I136").Value = 2251.125
$ws.Range("K136").Value = 6753.375
$ws.Range("M136").Value = -4203.375

$ws.Range("H141").Value = 375997.7
$ws.Range("J141").Value = 375997.7
$ws.Range("L141").Value = 375997.7
$ws.Range("N141").Value = -386357.7

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 313240.5
$ws.Range("J37").Value = 313240.5
$ws.Range("L37").Value = 939721.5
$ws.Range("N37").Value = -939945.5

$ws.Range("H41").Value = 47620280
$ws.Range("J41").Value = 55556740
$ws.Range("L41").Value = 166670220
$ws.Range("N41").Value = -166670896

$ws.Range("H129").Value = 2493.25
$ws.Range("J129").Value = 2493.25
$ws.Range("L129").Value = 7479.75
$ws.Range("N129").Value = -17479.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 37529.234
$ws.Range("I70").Value = 52910.74
$ws.Range("J70").Value = 5367.909
$ws.Range("K70").Value = 52910.74
$ws.Range("L70").Value = 5367.909
$ws.Range("M70").Value = -52640.74
$ws.Range("N70").Value = -5907.909

$ws.Range("H73").Value = 37529.234
$ws.Range("I73").Value = 52910.74
$ws.Range("J73").Value = 5367.909
$ws.Range("K73").Value = 52910.74
$ws.Range("L73").Value = 5367.909
$ws.Range("M73").Value = -51974.74
$ws.Range("N73").Value = -7239.909

$ws.Range("H80").Value = 25059486
$ws.Range("I80").Value = 78289.336
$ws.Range("J80").Value = 100003080
$ws.Range("K80").Value = 78289.336
$ws.Range("L80").Value = 100003080
$ws.Range("M80").Value = -77291.336
$ws.Range("N80").Value = -100005076

$ws.Range("H83").Value = 25059486
$ws.Range("I83").Value = 78289.336
$ws.Range("J83").Value = 100003080
$ws.Range("K83").Value = 391446.68
$ws.Range("L83").Value = 500015400
$ws.Range("M83").Value = -386454.68
$ws.Range("N83").Value = -500025384

$ws.Range("H102").Value = 2617.5454
$ws.Range("I102").Value = 2503.4285
$ws.Range("J102").Value = 5014
$ws.Range("K102").Value = 2503.4285
$ws.Range("L102").Value = 5014
$ws.Range("M102").Value = -881.4285
$ws.Range("N102").Value = -8258

$ws.Range("H132").Value = 8848.880999999999
$ws.Range("I132").Value = 7361.8
$ws.Range("K132").Value = 22085.4
$ws.Range("M132").Value = -19555.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5284.6924
$ws.Range("I7").Value = 2801.75
$ws.Range("J7").Value = 6388.222
$ws.Range("K7").Value = 2801.75
$ws.Range("L7").Value = 6388.222
$ws.Range("M7").Value = -2689.75
$ws.Range("N7").Value = -6612.222

$ws.Range("H22").Value = 1575.75
$ws.Range("I22").Value = 1533.6666
$ws.Range("J22").Value = 1702
$ws.Range("K22").Value = 1533.6666
$ws.Range("L22").Value = 1702
$ws.Range("M22").Value = -1238.6666
$ws.Range("N22").Value = -2292

$ws.Range("H27").Value = 1575.75
$ws.Range("I27").Value = 1533.6666
$ws.Range("J27").Value = 1702
$ws.Range("K27").Value = 1533.6666
$ws.Range("L27").Value = 1702
$ws.Range("M27").Value = -1426.6666
$ws.Range("N27").Value = -1916

$ws.Range("H126").Value = 5284.6924
$ws.Range("I126").Value = 2801.75
$ws.Range("J126").Value = 6388.222
$ws.Range("K126").Value = 8405.25
$ws.Range("L126").Value = 19164.666
$ws.Range("M126").Value = -5935.25
$ws.Range("N126").Value = -24104.666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 1116.4
$ws.Range("I17").Value = 1380.5714
$ws.Range("J17").Value = 500
$ws.Range("K17").Value = 1380.5714
$ws.Range("L17").Value = 500
$ws.Range("M17").Value = -1208.5714
$ws.Range("N17").Value = -844

$ws.Range("H29").Value = 7252.5
$ws.Range("I29").Value = 7252.5
$ws.Range("K29").Value = 7252.5
$ws.Range("M29").Value = -6962.5

$ws.Range("I81").Value = 5000
$ws.Range("K81").Value = 10000
$ws.Range("M81").Value = -8939

$ws.Range("I84").Value = 5000
$ws.Range("K84").Value = 50000
$ws.Range("M84").Value = -44696

$ws.Range("H107").Value = 1875.7142
$ws.Range("J107").Value = 1980
$ws.Range("L107").Value = 5940
$ws.Range("N107").Value = -9780

$ws.Range("H126").Value = 3097.5
$ws.Range("I126").Value = 3097.5
$ws.Range("K126").Value = 9292.5
$ws.Range("M126").Value = -6822.5

$ws.Range("H132").Value = 3718.5
$ws.Range("I132").Value = 2498.0732
$ws.Range("K132").Value = 7494.219599999999
$ws.Range("M132").Value = -4964.219599999999

$ws.Range("H136").Value = 4160.857
$ws.Range("J136").Value = 4258.9473
$ws.Range("L136").Value = 12776.8419
$ws.Range("N136").Value = -17876.8419
